# RMA Complete Flow(Issue Credit) template update
# "Misc Change-CPQ CRUD operations"
#
# The "RMA Details Maintenance Grid" sheet holds 3 sample RMA detail rows
# (rows 2-4). The Sales-Order-Line / Shipper-Line / Id columns (E, F, J)
# referenced an older RMA record ("RMA-CL8I-*"). This change swaps them to
# a freshly generated RMA record ("RMA-U8BT-*") with new Salesforce Ids,
# matching the pattern already used by every previous RMA batch in this
# grid (CL8I, TAFS, C084, ... -> U8BT).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Row 2 (Lot-track product line)
$ws.Range("E2").Value = "RMA-U8BT-001"
$ws.Range("F2").Value = "a7s5f000000xK6VAAU"
$ws.Range("J2").Value = "RMA-U8BT-1-1"

# Row 3 (Stock product, no-track line)
$ws.Range("E3").Value = "RMA-U8BT-002"
$ws.Range("F3").Value = "a7s5f000000xK6WAAU"
$ws.Range("J3").Value = "RMA-U8BT-1-2"

# Row 4 (Serial-track product line)
$ws.Range("E4").Value = "RMA-U8BT-003"
$ws.Range("F4").Value = "a7s5f000000xK6XAAU"
$ws.Range("J4").Value = "RMA-U8BT-1-3"

# Columns F (Shipper Line) and J (Id) are best-fit/auto-fit; the new
# values change their natural widths, so nudge the stored column widths
# to track the new content (closest achievable snap values).
$ws.Columns.Item(6).ColumnWidth = 13.666666666666666
$ws.Columns.Item(10).ColumnWidth = 20.833333333333336
